# Rubrica de Evaluacion - apply the commit's edits via Excel COM interop

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# ---------------------------------------------------------------------------
# 1) Update the four "Observacion" comments (column I) for each student
# ---------------------------------------------------------------------------
$ws.Range("I4").Value = "Se puso a pilas con la caja negra, echandole mayor enfoque a eso"
$ws.Range("I5").Value = "Se puso a pilas con la codificación, aunque aguantó para desvelarse, y no aguantó por un error que cometió"
$ws.Range("I6").Value = "Se reunió con Byron para hacer lo mejor, necesitó descansar en la madrugada, y lamentablemente cometió un fallo hasta que se le aruinara el Visual"
$ws.Range("I7").Value = "Hizo lo que pudo en el manual"

# ---------------------------------------------------------------------------
# 2) Update the scoring grid (columns D:G, rows 4-7). H holds =SUM(D:G) already.
# ---------------------------------------------------------------------------
$ws.Range("D4:G4").Value = 0.25
$ws.Range("D5:G5").Value = 0.25
$ws.Range("D6:F6").Value = 0.25
$ws.Range("G6").Value = 0
$ws.Range("D7:G7").Value = 0.25

# ---------------------------------------------------------------------------
# 3) Row heights
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 90
$ws.Rows.Item(6).RowHeight = 120
$ws.Rows.Item(7).RowHeight = 30

# ---------------------------------------------------------------------------
# 4) Selection / view state: select E6 (also clears the old topLeftCell scroll)
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Range("E6").Select() | Out-Null
